$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" value (row 8, column B) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2022-09-02T15:43:08-05:00"

# --- "Include from NIH Common Data " sheet: replace the per-concept table
#     (Concept/CADSR.../..."Overlap cGVHD") with a single "Codes"/"All codes"
#     pair, keeping the trailing blank row and the "System URI" row. ---
$ws = $wb.Worksheets.Item("Include from NIH Common Data ")

# Remove the four middle rows that enumerated individual CADSR concepts
# (old rows 3-6: CADSR:2991794.. through CADSR:2991795/Overlap cGVHD),
# shifting the blank row + System URI row up to become rows 3-4.
$ws.Range("A3:A6").EntireRow.Delete()

# Row 1 used to be "Concept" / "CADSR:4722619" -> now just "Codes" (no B cell).
$ws.Range("A1").Value = "Codes"
$ws.Range("B1").Clear()

# Row 2 used to be "CADSR:4722619" / "Acute GVHD" -> now just "All codes" (no B cell).
$ws.Range("A2").Value = "All codes"
$ws.Range("B2").Clear()
